# Meeting-diary update: restyle row 11, append rows 12-14 on "Part B" with
# the new meeting-log entries (report readiness, heatmap fixes, presentation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part B")

# --- Free up the old "fill/border" alignment styling used by row 11 so it
#     reads via the plain center/center (and wrap) styling that already
#     exists elsewhere on the sheet (same visual result: no fill or border
#     was ever actually painted, since fillId=0/borderId=0). ---
$ws.Range("B11:C11").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B11:C11").VerticalAlignment = -4108     # xlCenter

$ws.Range("E11").HorizontalAlignment = -4108       # xlCenter
$ws.Range("E11").VerticalAlignment = -4108         # xlCenter
$ws.Range("E11").WrapText = $true

# --- Values for the three new meeting rows ---

# Row 12 : 25 Oct 2023, 19:00-21:00
$ws.Range("A12").Value = 45224
$ws.Range("B12").Value = 0.79166666666666663
$ws.Range("C12").Value = 0.875
$ws.Range("D12").Value = "All"
$ws.Range("E12").Value = "Working on report readiness, removing redundant plots, finalised temporal and eco tourism analysis and improved overall report coherency"

# Row 13 : 26 Oct 2023, 15:00-18:00
$ws.Range("A13").Value = 45225
$ws.Range("B13").Value = 0.625
$ws.Range("C13").Value = 0.75
$ws.Range("D13").Value = "All"
$ws.Range("E13").Value = "Corected aspects of the heatmap, made all the plot themes uniform, allocated sections for the presentation"

# Row 14 : 26 Oct 2023, 21:00-22:00
$ws.Range("A14").Value = 45225
$ws.Range("B14").Value = 0.875
$ws.Range("C14").Value = 0.91666666666666663
$ws.Range("D14").Value = "All + Paired presentation team"
$ws.Range("E14").Value = "Presented report, participated in QnA"

# --- Number formats: format one cell, then propagate the exact same style
#     record to its siblings via copy/paste-format so the three rows share
#     a single style-table entry instead of three near-duplicate ones. ---
$ws.Range("A12").NumberFormat = "mm-dd-yy"
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B12").NumberFormat = "h:mm AM/PM"
$ws.Range("B12").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C12:C14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Alignment for the text columns (matches existing sheet styling) ---
$ws.Range("D12:D13").HorizontalAlignment = -4108
$ws.Range("D12:D13").VerticalAlignment = -4108

$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("D14").VerticalAlignment = -4108
$ws.Range("D14").WrapText = $true

$ws.Range("E12:E14").HorizontalAlignment = -4108
$ws.Range("E12:E14").VerticalAlignment = -4108
$ws.Range("E12:E14").WrapText = $true

# --- Row heights to fit the wrapped discussion text ---
$ws.Rows.Item(12).RowHeight = 153
$ws.Rows.Item(13).RowHeight = 102
$ws.Rows.Item(14).RowHeight = 34

# --- View: scroll so row 11 is at the top, selection sits on E16 ---
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("E16").Select()

Write-Output "meeting diary updated"
